$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos -> professor name leaked into this row ---
$ws.Range("B10").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C10").Value = "11079086 - Herlandí de Souza Andrade"

# --- Row 13: used to be a headerless row holding the professor name; ---
# --- it becomes "Programa resumido:" / "Semestral"                   ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# --- Row 14: "Programa resumido:" -> "Short syllabus:" ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "To be defined according to the scheduled topic"
$ws.Range("C14").Value = "To be defined according to the scheduled topic"

# --- Row 15: "Short syllabus:" -> "Programa:" ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2021"
$ws.Range("C15").Value = "01/01/2021"
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16: "Programa:" -> "Syllabus:" (height stays 120) ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "The content of this optional course will be according to the topic to be programmed, and should address complementary subjects to the regular content of the undergraduate course."
$ws.Range("C16").Value = "The content of this optional course will be according to the topic to be programmed, and should address complementary subjects to the regular content of the undergraduate course."

# --- Row 17: "Syllabus:" -> "Avaliação:" (no B/C, default height) ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = ""
$ws.Rows.Item(17).EntireRow.AutoFit()

# --- Row 18: "Avaliação:" -> "Método:" with the professor name again ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C18").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19: "Método:" -> "Critério:" (B/C text unchanged) ---
$ws.Range("A19").Value = "Critério:"

# --- Row 20: "Critério:" -> "Norma de recuperação:" (B/C text unchanged) ---
$ws.Range("A20").Value = "Norma de recuperação:"

# --- Row 21: "Norma de recuperação:" -> "Bibliografia:" (B/C text unchanged) ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows.Item(21).RowHeight = 120

# --- Row 22 is dropped entirely (it used to hold "Bibliografia:" / the old bibliography text) ---
$ws.Rows.Item(22).Delete()
